$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.672.21"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.689.61"
$ws.Range("E3").Value = "  -0.91%  "
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.33"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3940"
$ws.Range("E7").Value = "  -0.83%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4049"
$ws.Range("E8").Value = "  -0.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.486"
$ws.Range("E9").Value = "  -1.69%  "
$ws.Range("E10").Value = "  +0.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.53"
$ws.Range("E11").Value = "  -1.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08827"
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.237"
$ws.Range("E13").Value = "  -1.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.48"
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.033"
$ws.Range("E15").Value = "  +6.00%  "
$ws.Range("E16").Value = "  -0.96%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.691.24"
$ws.Range("E17").Value = "  -0.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "99.52"
$ws.Range("E18").Value = "  -1.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07011"
$ws.Range("E19").Value = "  -1.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.49"
$ws.Range("E20").Value = "  -0.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.987"
$ws.Range("E21").Value = "  +2.99%  "
$ws.Range("E22").Value = "  +0.33%  "
$ws.Range("E23").Value = "  +0.83%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.657.71"
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.304"
$ws.Range("E25").Value = "  +10.00%  "
$ws.Range("E26").Value = "  +1.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.75"
$ws.Range("E27").Value = "  +0.98%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.53"
$ws.Range("E28").Value = "  +1.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "135.37"
$ws.Range("E29").Value = "  +1.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.177"
$ws.Range("E30").Value = "  +0.76%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.623"
$ws.Range("E31").Value = "  +2.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.878.93"
$ws.Range("E32").Value = "  -0.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.057"
$ws.Range("E33").Value = "  -2.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08532"
$ws.Range("E34").Value = "  -2.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.096"
$ws.Range("E35").Value = "  -2.57%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.20"
$ws.Range("E36").Value = "  +0.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2732"
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.887"
$ws.Range("E38").Value = "  -3.80%  "
$ws.Range("E39").Value = "  -2.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09182"
$ws.Range("E40").Value = "  +1.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02719"
$ws.Range("E41").Value = "  -2.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.465"
$ws.Range("E42").Value = "  -1.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7619"
$ws.Range("E43").Value = "  -1.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.11"
$ws.Range("E44").Value = "  +2.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.583"
$ws.Range("E45").Value = "  +4.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.7127"
$ws.Range("E46").Value = "  -1.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.209"
$ws.Range("E47").Value = "  +0.89%  "
$ws.Range("E48").Value = "  +0.30%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "139.83"
$ws.Range("E49").Value = "  -1.21%  "
$ws.Range("B50").Value = "Flow"
$ws.Range("C50").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.317"
$ws.Range("E50").Value = "  +1.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07968"
$ws.Range("E51").Value = "  -0.48%  "
